$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = "16.11.2015."
$ws.Range("C9").Value = "index.php"
$ws.Range("D9").Value = "Pēc autorizācijas vajag atvērt atbilstošo lapu katrai lomai."
$ws.Range("E9").Value = "J"

$ws.Rows.Item(9).RowHeight = 30

$ws.Range("G8").Select()
